# Update column O ("player_id_y") values on Sheet1 to match the refreshed
# CSV data downloaded from kaggle.com. Only the player_id_y values for a
# subset of rows changed; every other cell stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2;  Old = 5025; New = 5032 },
    @{ Row = 3;  Old = 3251; New = 3257 },
    @{ Row = 4;  Old = 2223; New = 2227 },
    @{ Row = 5;  Old = 2936; New = 2941 },
    @{ Row = 6;  Old = 609;  New = 610  },
    @{ Row = 7;  Old = 3966; New = 3972 },
    @{ Row = 8;  Old = 5061; New = 5068 },
    @{ Row = 9;  Old = 1947; New = 1951 },
    @{ Row = 10; Old = 718;  New = 719  },
    @{ Row = 11; Old = 2638; New = 2643 },
    @{ Row = 12; Old = 1970; New = 1974 },
    @{ Row = 13; Old = 3994; New = 4000 },
    @{ Row = 14; Old = 4452; New = 4459 },
    @{ Row = 15; Old = 3104; New = 3109 },
    @{ Row = 17; Old = 1142; New = 1144 },
    @{ Row = 18; Old = 2091; New = 2095 },
    @{ Row = 19; Old = 3772; New = 3778 },
    @{ Row = 20; Old = 582;  New = 583  },
    @{ Row = 22; Old = 5265; New = 5272 },
    @{ Row = 23; Old = 1951; New = 1955 },
    @{ Row = 24; Old = 2516; New = 2521 },
    @{ Row = 25; Old = 4941; New = 4948 },
    @{ Row = 26; Old = 2321; New = 2325 },
    @{ Row = 27; Old = 1173; New = 1175 },
    @{ Row = 29; Old = 2174; New = 2178 },
    @{ Row = 30; Old = 1945; New = 1949 },
    @{ Row = 31; Old = 1122; New = 1124 },
    @{ Row = 32; Old = 1461; New = 1464 },
    @{ Row = 33; Old = 3404; New = 3410 },
    @{ Row = 34; Old = 4038; New = 4044 },
    @{ Row = 35; Old = 4525; New = 4532 },
    @{ Row = 36; Old = 4411; New = 4418 },
    @{ Row = 37; Old = 3097; New = 3102 },
    @{ Row = 38; Old = 3088; New = 3093 },
    @{ Row = 39; Old = 1085; New = 1087 }
)

$colIndex = 15  # Column O

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $colIndex)
    $cell.Value = $u.New
}
